$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.080.58"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "2.469.43"
$ws.Range("E3").Value = "  -2.53%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.95%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.515"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("D9").Value = "2.469.54"
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("E10").Value = "  -2.41%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.67%  "
$ws.Range("E13").Value = "  -3.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.37%  "
$ws.Range("D16").Value = "67.007.27"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("E17").Value = "  -4.37%  "
$ws.Range("D18").Value = "2.472.43"
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("E19").Value = "  -5.76%  "
$ws.Range("E20").Value = "  -3.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "354.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.94%  "
$ws.Range("E25").Value = "  -7.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.39%  "
$ws.Range("E28").Value = "  -6.70%  "
$ws.Range("D29").Value = "2.566.50"
$ws.Range("E29").Value = "  -3.31%  "
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "518.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.35%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0907"
$ws.Range("E31").Value = "  -6.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.81%  "
$ws.Range("E33").Value = "  -5.65%  "
$ws.Range("E34").Value = "  -5.88%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  -6.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.50"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.44%  "
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("E39").Value = "  -3.57%  "
$ws.Range("E40").Value = "  -5.72%  "
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.59%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.49%  "
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.327"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.63%  "
$ws.Range("E48").Value = "  -6.74%  "
$ws.Range("E49").Value = "  -6.77%  "
$ws.Range("D50").Value = "0.0₆0255"
$ws.Range("E50").Value = "  -11.38%  "
$ws.Range("E51").Value = "  -7.21%  "
